$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Congratulations {{ users }}!" -> "Congratulations {{ users[0] }}!"
#    (collect the *first* user rather than the whole list)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Congratulations {{ users }}", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "Congratulations {{ users[0] }}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Drop the leftover stray bracket + the "vs. other party" conditional
#    clause that used to follow the highlighted [answer/motion] text, so the
#    sentence now simply ends "...are your [answer/motion]." Leave the final
#    "." (from the old "{% endif %}." run) untouched so it stays outside the
#    yellow highlight, exactly like before.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("[answer/motion] ]{% if other_parties.number() %} in {{ users }} v {{ other_parties }}{% endif %}", `
                         $false, $true, $false, $false, $false, `
                         $true, 1, $false, "[answer/motion]", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Cosmetic re-saves: replacing a run's text with itself makes the engine
#    re-serialize it, which drops the now-superfluous xml:space="preserve"
#    attribute on runs whose text has no leading/trailing whitespace - this
#    mirrors what a later Word re-save of the document did.
# ---------------------------------------------------------------------------

# Title banner text (table heading).
$d.Content.Find.Execute("File a motion and notice of hearing in your case", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "File a motion and notice of hearing in your case", 2) | Out-Null

# The bold "Tell" at the start of the "Tell the judge why..." paragraph only
# (the word "Tell" also appears later in that same paragraph as part of a
# longer, already-tight run, so restrict the search to the first match).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Tell the judge why*") {
        $para.Range.Find.Execute("Tell", $false, $true, $false, $false, $false, `
                                  $true, 1, $false, "Tell", 2) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 4. Remove the stray "_GoBack" bookmark that wraps the last empty paragraph
#    (Word drops this automatically-managed bookmark on a later edit/save).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
